$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2343023818649499
$ws.Range("C2").Value = 0.04130090163999256
$ws.Range("D2").Value = 0.02907139404816661
$ws.Range("F2").Value = 0.6981695318995662
$ws.Range("G2").Value = 0.002425462473374959
$ws.Range("K2").Value = 0.1991907139332341
$ws.Range("M2").Value = 0.8862292840141777
$ws.Range("N2").Value = 1.474153485319732
$ws.Range("O2").Value = 2.391647609380868
$ws.Range("B3").Value = 0.2082200658922204
$ws.Range("C3").Value = 0.0394779188220582
$ws.Range("D3").Value = 0.02699843040836925
$ws.Range("F3").Value = 0.695350594744184
$ws.Range("G3").Value = 0.002427727375901607
$ws.Range("K3").Value = 0.1739934342470377
$ws.Range("M3").Value = 0.7891297409906883
$ws.Range("N3").Value = 1.491270294517424
$ws.Range("O3").Value = 2.394980879541748
$ws.Range("B4").Value = 0.1922216416344327
$ws.Range("C4").Value = 0.03834714272699813
$ws.Range("D4").Value = 0.02571295495988579
$ws.Range("F4").Value = 0.6940259653743297
$ws.Range("G4").Value = 0.002429191739130007
$ws.Range("K4").Value = 0.1584864544024356
$ws.Range("M4").Value = 0.7300310314882807
$ws.Range("N4").Value = 1.502310626141941
$ws.Range("O4").Value = 2.398329195652224
$ws.Range("B5").Value = 0.1857066839226604
$ws.Range("C5").Value = 0.03788350752345337
$ws.Range("D5").Value = 0.02518596188409816
$ws.Range("F5").Value = 0.6935882912189015
$ws.Range("G5").Value = 0.002429807071865565
$ws.Range("K5").Value = 0.1521586904472514
$ws.Range("M5").Value = 0.7060729073862433
$ws.Range("N5").Value = 1.506942955706435
$ws.Range("O5").Value = 2.400020814841142
$ws.Range("B6").Value = 0.1846251670241656
$ws.Range("C6").Value = 0.03780635139562349
$ws.Range("D6").Value = 0.02509826574573282
$ws.Range("F6").Value = 0.693521783137129
$ws.Range("G6").Value = 0.002429910372072297
$ws.Range("K6").Value = 0.1511074671442998
$ws.Range("M6").Value = 0.7021020636914415
$ws.Range("N6").Value = 1.507720198777438
$ws.Range("O6").Value = 2.400321461467968
$ws.Range("B7").Value = 0.1921337596967874
$ws.Range("C7").Value = 0.0383409013984064
$ws.Range("D7").Value = 0.02570586046186918
$ws.Range("F7").Value = 0.6940196492865169
$ws.Range("G7").Value = 0.002429199962442388
$ws.Range("K7").Value = 0.1584011499683697
$ws.Range("M7").Value = 0.7297074250801927
$ws.Range("N7").Value = 1.502372559764542
$ws.Range("O7").Value = 2.398350685017888
$ws.Range("B8").Value = 0.225306094454794
$ws.Range("C8").Value = 0.04067474238141244
$ws.Range("D8").Value = 0.02835928511102992
$ws.Range("F8").Value = 0.6971132374856523
$ws.Range("G8").Value = 0.002426228150156407
$ws.Range("K8").Value = 0.1905103711565204
$ws.Range("M8").Value = 0.8526376725037181
$ws.Range("N8").Value = 1.47994518598145
$ws.Range("O8").Value = 2.392526653463335
$ws.Range("B9").Value = 0.2904690256945344
$ws.Range("C9").Value = 0.04515863934734909
$ws.Range("D9").Value = 0.03346085499073581
$ws.Range("F9").Value = 0.706405192592328
$ws.Range("G9").Value = 0.002420982537670299
$ws.Range("K9").Value = 0.2531761557849848
$ws.Range("M9").Value = 1.098116391001597
$ws.Range("N9").Value = 1.440178767550872
$ws.Range("O9").Value = 2.391444174997531
$ws.Range("B10").Value = 0.3383948406939794
$ws.Range("C10").Value = 0.04839419070436435
$ws.Range("D10").Value = 0.03714552092529999
$ws.Range("F10").Value = 0.715203545885899
$ws.Range("G10").Value = 0.002417479651504852
$ws.Range("K10").Value = 0.299015537500253
$ws.Range("M10").Value = 1.281570784205456
$ws.Range("N10").Value = 1.413535589316993
$ws.Range("O10").Value = 2.396969749277531
$ws.Range("B11").Value = 0.3602050338549532
$ws.Range("C11").Value = 0.04985290802301279
$ws.Range("D11").Value = 0.03880770952491019
$ws.Range("F11").Value = 0.7196355188453936
$ws.Range("G11").Value = 0.002415961516911196
$ws.Range("K11").Value = 0.3198217799314875
$ws.Range("M11").Value = 1.365794437404716
$ws.Range("N11").Value = 1.401974655637243
$ws.Range("O11").Value = 2.400860275768821
$ws.Range("B12").Value = 0.3684648076774693
$ws.Range("C12").Value = 0.05040335026421872
$ws.Range("D12").Value = 0.03943509469959139
$ws.Range("F12").Value = 0.7213756280023347
$ws.Range("G12").Value = 0.002415397411740208
$ws.Range("K12").Value = 0.3276934987398761
$ws.Range("M12").Value = 1.397805737197814
$ws.Range("N12").Value = 1.397677365554788
$ws.Range("O12").Value = 2.402531801092579
$ws.Range("B13").Value = 0.366685893154056
$ws.Range("C13").Value = 0.05028488988169499
$ws.Range("D13").Value = 0.03930006784861462
$ws.Range("F13").Value = 0.7209981142628834
$ws.Range("G13").Value = 0.002415518423224113
$ws.Range("K13").Value = 0.3259985092606712
$ws.Range("M13").Value = 1.390906186071916
$ws.Range("N13").Value = 1.398599277006991
$ws.Range("O13").Value = 2.402162985886775
$ws.Range("B14").Value = 0.3608845586860241
$ws.Range("C14").Value = 0.04989823246005187
$ws.Range("D14").Value = 0.03885936617901109
$ws.Range("F14").Value = 0.7197774394911249
$ws.Range("G14").Value = 0.002415914891939857
$ws.Range("K14").Value = 0.3204695371685204
$ws.Range("M14").Value = 1.368425628094101
$ws.Range("N14").Value = 1.401619498945109
$ws.Range("O14").Value = 2.400993818037193
$ws.Range("B15").Value = 0.3573311537883512
$ws.Range("C15").Value = 0.04966113914963444
$ws.Range("D15").Value = 0.03858915561828979
$ws.Range("F15").Value = 0.7190377927851941
$ws.Range("G15").Value = 0.002416159142615052
$ws.Range("K15").Value = 0.3170819358718404
$ws.Range("M15").Value = 1.35467116690846
$ws.Range("N15").Value = 1.403479971841428
$ws.Range("O15").Value = 2.400303497275701
$ws.Range("B16").Value = 0.3369696219880041
$ws.Range("C16").Value = 0.04829859108408385
$ws.Range("D16").Value = 0.03703660828875144
$ws.Range("F16").Value = 0.7149225538571997
$ws.Range("G16").Value = 0.002417580376643075
$ws.Range("K16").Value = 0.2976548264104224
$ws.Range("M16").Value = 1.276082645873601
$ws.Range("N16").Value = 1.414302390548308
$ws.Range("O16").Value = 2.396743224146945
$ws.Range("B17").Value = 0.3244803064532107
$ws.Range("C17").Value = 0.0474593072250542
$ws.Range("D17").Value = 0.03608056159756501
$ws.Range("F17").Value = 0.7125080415463216
$ws.Range("G17").Value = 0.002418471516202184
$ws.Range("K17").Value = 0.2857247094698039
$ws.Range("M17").Value = 1.228073105215969
$ws.Range("N17").Value = 1.421084936904542
$ws.Range("O17").Value = 2.394911962748125
$ws.Range("B18").Value = 0.3172976121268505
$ws.Range("C18").Value = 0.04697533856008107
$ws.Range("D18").Value = 0.0355293548650053
$ws.Range("F18").Value = 0.7111597070900686
$ws.Range("G18").Value = 0.002418991171244991
$ws.Range("K18").Value = 0.2788584911507712
$ws.Range("M18").Value = 1.200531271127048
$ws.Range("N18").Value = 1.425038700313113
$ws.Range("O18").Value = 2.393988255119865
$ws.Range("B19").Value = 0.3148658328182137
$ws.Range("C19").Value = 0.04681126456810603
$ws.Range("D19").Value = 0.03534250123775706
$ws.Range("F19").Value = 0.7107101265259814
$ws.Range("G19").Value = 0.002419168337969842
$ws.Range("K19").Value = 0.2765329788223596
$ws.Range("M19").Value = 1.191218232118686
$ws.Range("N19").Value = 1.426386409295318
$ws.Range("O19").Value = 2.393697751883849
$ws.Range("B20").Value = 0.3258097331330703
$ws.Range("C20").Value = 0.04754877859236473
$ws.Range("D20").Value = 0.03618247066967228
$ws.Range("F20").Value = 0.7127608858998471
$ws.Range("G20").Value = 0.002418375918943051
$ws.Range("K20").Value = 0.2869951425498698
$ws.Range("M20").Value = 1.233176302545786
$ws.Range("N20").Value = 1.420357475822975
$ws.Range("O20").Value = 2.395093490114164
$ws.Range("B21").Value = 0.362588535463118
$ws.Range("C21").Value = 0.05001185623983417
$ws.Range("D21").Value = 0.03898886687567682
$ws.Range("F21").Value = 0.7201343031798615
$ws.Range("G21").Value = 0.002415798147334192
$ws.Range("K21").Value = 0.3220937283726073
$ws.Range("M21").Value = 1.375025463802231
$ws.Range("N21").Value = 1.400730198092294
$ws.Range("O21").Value = 2.40133184789417
$ws.Range("B22").Value = 0.3866296163721188
$ws.Range("C22").Value = 0.0516102745070981
$ws.Range("D22").Value = 0.04081104451772433
$ws.Range("F22").Value = 0.7253135979747611
$ws.Range("G22").Value = 0.00241417622995499
$ws.Range("K22").Value = 0.3449907572517645
$ws.Range("M22").Value = 1.468421103747616
$ws.Range("N22").Value = 1.388372430250197
$ws.Range("O22").Value = 2.406564756916083
$ws.Range("B23").Value = 0.373798236090721
$ws.Range("C23").Value = 0.05075822365282079
$ws.Range("D23").Value = 0.03983962156905818
$ws.Range("F23").Value = 0.7225163238670973
$ws.Range("G23").Value = 0.00241503614926053
$ws.Range("K23").Value = 0.3327741811473857
$ws.Range("M23").Value = 1.418508689540928
$ws.Range("N23").Value = 1.394924958572208
$ws.Range("O23").Value = 2.40366601093919
$ws.Range("B24").Value = 0.3252087068328819
$ws.Range("C24").Value = 0.04750833310874469
$ws.Range("D24").Value = 0.03613640244326177
$ws.Range("F24").Value = 0.7126464510004169
$ws.Range("G24").Value = 0.002418419115756419
$ws.Range("K24").Value = 0.2864208028759663
$ws.Range("M24").Value = 1.230868962261184
$ws.Range("N24").Value = 1.420686191587372
$ws.Range("O24").Value = 2.395011019466693
$ws.Range("B25").Value = 0.2728305581805728
$ws.Range("C25").Value = 0.04395578605002015
$ws.Range("D25").Value = 0.03209178710967109
$ws.Range("F25").Value = 0.7035456281412351
$ws.Range("G25").Value = 0.002422339690695785
$ws.Range("K25").Value = 0.2362574199590881
$ws.Range("M25").Value = 1.031193532266926
$ws.Range("N25").Value = 1.450485094428196
$ws.Range("O25").Value = 2.390628337834158

Write-Host "Updated cells"